# "architecture nav for intro"
# The two "KNX TP" device boxes on the architecture slide are relabelled to "KNX IP".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

foreach ($shp in $s.Shapes) {
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        $tr = $shp.TextFrame.TextRange
        foreach ($para in $tr.Paragraphs()) {
            if ($para.Text.Trim() -eq "KNX TP") {
                $para.Runs(1).Text = "KNX IP"
            }
        }
    }
}
